$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data cell in this sheet is stored as literal text (inline strings),
# even when the text looks numeric (e.g. "330.32", "-0.81%", the "0"/"2" flag
# in column G). Plain `.Value = "..."` assignment lets Excel reinterpret a
# numeric-looking string as a real number/percentage, so each touched cell is
# briefly marked Text (NumberFormat "@") while the literal value is written,
# then ClearFormats() drops the temporary number-format mark again so the
# cell keeps its original (unformatted) appearance with the value stored as text.
function Set-TextValue([string]$cell, [string]$value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "333.98"
Set-TextValue "E2" "0.43%"
Set-TextValue "G2" "2"

# Row 3
Set-TextValue "D3" "42.80"
Set-TextValue "E3" "2.60%"
Set-TextValue "G3" "2"

# Row 4
Set-TextValue "D4" "5.651"
Set-TextValue "E4" "-0.57%"
Set-TextValue "G4" "2"

# Row 5
Set-TextValue "D5" "0.08292"
Set-TextValue "E5" "0.95%"
Set-TextValue "G5" "2"

# Row 6
Set-TextValue "D6" "8.800"
Set-TextValue "E6" "0.31%"
Set-TextValue "G6" "2"

# Row 7
Set-TextValue "D7" "4.502"
Set-TextValue "E7" "-0.89%"
Set-TextValue "G7" "2"

# Row 8
Set-TextValue "D8" "1.963"
Set-TextValue "E8" "-3.83%"
Set-TextValue "G8" "2"

# Row 9
Set-TextValue "D9" "2.899"
Set-TextValue "E9" "-0.88%"
Set-TextValue "G9" "2"

# Row 10
Set-TextValue "D10" "0.9211"
Set-TextValue "E10" "-0.31%"
Set-TextValue "G10" "2"

# Row 11
Set-TextValue "D11" "0.1245"
Set-TextValue "E11" "-1.24%"
Set-TextValue "G11" "2"

# Row 12
Set-TextValue "D12" "0.1944"
Set-TextValue "E12" "-0.66%"
Set-TextValue "G12" "2"

# Row 13
Set-TextValue "D13" "0.09330"
Set-TextValue "E13" "-2.06%"
Set-TextValue "G13" "2"

# Row 14
Set-TextValue "D14" "0.03978"
Set-TextValue "E14" "6.95%"
Set-TextValue "G14" "2"

# Row 15
Set-TextValue "D15" "0.1066"
Set-TextValue "E15" "0.94%"
Set-TextValue "G15" "2"

# Row 16
Set-TextValue "D16" "0.001307"
Set-TextValue "E16" "0.56%"
Set-TextValue "G16" "2"

# Row 17
Set-TextValue "D17" "0.005958"
Set-TextValue "E17" "-2.92%"
Set-TextValue "G17" "2"

# Row 18
Set-TextValue "D18" "3.496"
Set-TextValue "E18" "2.75%"
Set-TextValue "G18" "2"

# Row 19
Set-TextValue "E19" "0.74%"
Set-TextValue "G19" "2"

# Row 20
Set-TextValue "D20" "9.021"
Set-TextValue "E20" "9.42%"
Set-TextValue "G20" "2"

# Row 21
Set-TextValue "D21" "0.1373"
Set-TextValue "E21" "-2.21%"
Set-TextValue "G21" "2"

# Row 22
Set-TextValue "D22" "0.2577"
Set-TextValue "E22" "-2.87%"
Set-TextValue "G22" "2"

# Row 23
Set-TextValue "D23" "0.04422"
Set-TextValue "E23" "-0.32%"
Set-TextValue "G23" "2"

# Row 24
Set-TextValue "D24" "0.001241"
Set-TextValue "E24" "-2.55%"
Set-TextValue "G24" "2"

# Row 25
Set-TextValue "D25" "0.004448"
Set-TextValue "E25" "3.11%"
Set-TextValue "G25" "2"

# Row 26
Set-TextValue "D26" "0.0001193"
Set-TextValue "E26" "-0.75%"
Set-TextValue "G26" "2"

# Row 27
Set-TextValue "D27" "0.0003999"
Set-TextValue "E27" "0.16%"
Set-TextValue "G27" "2"

# Row 28
Set-TextValue "G28" "2"

# Row 29
Set-TextValue "G29" "2"

# Row 30
Set-TextValue "G30" "2"

# Row 31
Set-TextValue "G31" "2"

# Row 32
Set-TextValue "G32" "2"

# Row 33
Set-TextValue "G33" "2"

# Row 34
Set-TextValue "G34" "2"

# Row 35
Set-TextValue "G35" "2"

# Row 36
Set-TextValue "G36" "2"

# Row 37
Set-TextValue "G37" "2"

# Row 38
Set-TextValue "G38" "2"

# Row 39
Set-TextValue "D39" "0.02781"
Set-TextValue "E39" "-2.35%"
Set-TextValue "G39" "2"

# Row 40
Set-TextValue "D40" "0.05628"
Set-TextValue "E40" "2.94%"
Set-TextValue "G40" "2"

# Row 41
Set-TextValue "D41" "0.007910"
Set-TextValue "E41" "3.10%"
Set-TextValue "G41" "2"

# Row 42
Set-TextValue "D42" "0.1429"
Set-TextValue "E42" "0.69%"
Set-TextValue "G42" "2"

# Row 43
Set-TextValue "D43" "0.009096"
Set-TextValue "E43" "-3.62%"
Set-TextValue "G43" "2"

# Row 44
Set-TextValue "D44" "0.002180"
Set-TextValue "E44" "2.61%"
Set-TextValue "G44" "2"

# Row 45
Set-TextValue "D45" "0.01006"
Set-TextValue "E45" "-8.31%"
Set-TextValue "G45" "2"

# Row 46
Set-TextValue "D46" "0.00007318"
Set-TextValue "E46" "7.64%"
Set-TextValue "G46" "2"

# Row 47
Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "0.08%"
Set-TextValue "G47" "2"

# Row 48
Set-TextValue "D48" "0.003648"
Set-TextValue "E48" "13.47%"
Set-TextValue "G48" "2"

# Row 49
Set-TextValue "D49" "0.002284"
Set-TextValue "E49" "0.00%"
Set-TextValue "G49" "2"

# Row 50
Set-TextValue "D50" "0.00002106"
Set-TextValue "E50" "0.08%"
Set-TextValue "G50" "2"

# Row 51
Set-TextValue "D51" "0.0002006"
Set-TextValue "E51" "0.08%"
Set-TextValue "G51" "2"
